# Release abs09 2nd wave
#
# - Fix a wording typo in the German study description (sheet "study", cell C2):
#   "Absolventenkohorte" -> "Absolvent(inn)enkohorte" in the sentence about the
#   methodical innovation of the 2009 cohort.
# - Update the active window's view on the "study" sheet: scroll so column C is
#   the left-most visible column and move the selection to D8.

$wb = $excel.ActiveWorkbook

$study = $wb.Worksheets.Item("study")

# --- 1. Correct the German description text -----------------------------
$descCell = $study.Range("C2")
$oldText  = $descCell.Value2
$newText  = $oldText -replace `
    "methodische Neuerung der Absolventenkohorte 2009", `
    "methodische Neuerung der Absolvent(inn)enkohorte 2009"
$descCell.Value = $newText

# --- 2. Update the sheet view (scroll position + active selection) ------
$study.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$study.Range("D8").Select()
